$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 228
$ws1.Range("F4").Value = 143
$ws1.Range("F5").Value = 57
$ws1.Range("F6").Value = 3698
$ws1.Range("F7").Value = 219
$ws1.Range("F8").Value = 2494
$ws1.Range("F9").Value = 57
$ws1.Range("F10").Value = 2957
$ws1.Range("F12").Value = 522
$ws1.Range("F13").Value = 2252
$ws1.Range("F14").Value = 48
$ws1.Range("F16").Value = 37
$ws1.Range("F17").Value = 419
$ws1.Range("F19").Value = 180
$ws1.Range("F21").Value = 290
$ws1.Range("F22").Value = 302
$ws1.Range("F23").Value = 631
$ws1.Range("F24").Value = 1370
$ws1.Range("F25").Value = 34
$ws1.Range("F26").Value = 1283
$ws1.Range("F27").Value = 116
$ws1.Range("F30").Value = 4071
$ws1.Range("F31").Value = 3638
$ws1.Range("F34").Value = 1086
$ws1.Range("F35").Value = 438
$ws1.Range("F37").Value = 1299
$ws1.Range("F38").Value = 136
$ws1.Range("F42").Value = 52
$ws1.Range("F43").Value = 37

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 184
$ws2.Range("F6").Value = 18
$ws2.Range("F16").Value = 190

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1011
$ws3.Range("F4").Value = 137
$ws3.Range("F5").Value = 2186

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 228
$ws4.Range("F3").Value = 1011
$ws4.Range("F4").Value = 137
$ws4.Range("F7").Value = 143
$ws4.Range("F8").Value = 184
$ws4.Range("F9").Value = 57
$ws4.Range("F11").Value = 3698
$ws4.Range("F12").Value = 219
$ws4.Range("F13").Value = 2494
$ws4.Range("F14").Value = 57
$ws4.Range("F15").Value = 2957
$ws4.Range("F16").Value = 522
$ws4.Range("F17").Value = 2252
$ws4.Range("F18").Value = 48
$ws4.Range("F20").Value = 37
$ws4.Range("F21").Value = 419
$ws4.Range("F23").Value = 180
$ws4.Range("F25").Value = 302
$ws4.Range("F26").Value = 631
$ws4.Range("F27").Value = 1370
$ws4.Range("F28").Value = 34
$ws4.Range("F29").Value = 1283
$ws4.Range("F33").Value = 4071
$ws4.Range("F34").Value = 3638
$ws4.Range("F36").Value = 1086
$ws4.Range("F38").Value = 438
$ws4.Range("F43").Value = 1299
$ws4.Range("F44").Value = 136
$ws4.Range("F47").Value = 52
$ws4.Range("F48").Value = 37
$ws4.Range("F49").Value = 190
